$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before column N, shifting
# Late / heading / Outstanding columns one position to the right.
$ws.Columns("N").Insert()

# Excel copies the width of the column to the left (M) onto the
# newly inserted column.
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make "Repayment schedule" the active sheet/tab and select cell R10.
$ws.Activate()
$null = $ws.Range("R10").Select()
